$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rescale / re-bin the counts and add the extra time-bucket rows ---
#
# Before:
#   row1:        B1=0
#   row2: A2=0,  B2=306
#   row3: A3=1,  B3=163
#
# After:
#   row1:        B1=0                 (unchanged)
#   row2: A2=0,  B2=176               (value rescaled)
#   row3: A3=3,  B3=116               (new bucket)
#   row4: A4=1,  B4=89                (old row3 bucket, rescaled)
#   row5: A5=2,  B5=88                (new bucket)

# Capture the old row 3 contents before they get overwritten, so we can
# relocate them to row 5 further down.
$oldA3 = $ws.Range("A3").Value
$oldB3 = $ws.Range("B3").Value

# Apply the "A" column formatting (bold / bordered / centered, same as the
# other label cells in column A) to the two brand-new rows (3 and 4), and
# re-apply it to row 5 as well so everything stays visually consistent.
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Relocate the former row-3 bucket down to row 5 (value will be overwritten
# right after with the new, rescaled number).
$ws.Range("A5").Value = $oldA3
$ws.Range("B5").Value = $oldB3

# Rescale existing row 2 value.
$ws.Range("B2").Value = 176

# New row 3.
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 116

# Row 4 (rescaled old row-3 bucket).
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 89

# New row 5.
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 88
